$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the Gompertz "asymptote" parameter for dairy (column D) in D3.
# Downstream formulas in D9:D25 reference D3 and will recalc automatically.
$ws.Range("D3").Value = 0.018

# Update the last active selection to match the saved workbook state.
$ws.Range("D28:T28").Select()
